$d = $word.ActiveDocument

# 1. Remove the now-obsolete middle requirement run (its text plus the
#    trailing line-break character), leaving the other two runs intact.
$f = $d.Content
$f.Find.Execute(
    "LOM3229 -  Métodos Experimentais da Física II  (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delRange = $d.Range($f.Start, $f.End + 1)
$delRange.Delete()

# 2. Replace the first requirement line's text (run 1)
$d.Content.Find.Execute(
    "LOM3215 -  Física do Estado Sólido  (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LOM3241 -  Química de Materiais  (Requisito fraco)", 2)

# 3. Replace the remaining requirement line's text (run 3)
$d.Content.Find.Execute(
    "LOM3246 -  Técnicas de Caracterização de Materiais  (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LOB1021 -  Física IV  (Requisito fraco)", 2)
